# refatorando o consolidador para modelo ETL
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=22881; B="Ana Lívia Ramos";            C="Financeiro";             D="Problemas pessoais"; E=2; F=45084; G=9877.6}
    @{Row=3;  A=98313; B="Enzo Gabriel Moura";          C="P&D";                    D="Doença";              E=6; F=45090; G=9926.33}
    @{Row=4;  A=85482; B="Dr. Luiz Fernando Correia";   C="Financeiro";             D="Outros";              E=5; F=45090; G=8672.76}
    @{Row=5;  A=43820; B="Raquel Carvalho";             C="Marketing";              D="Problemas pessoais"; E=6; F=45100; G=10090.88}
    @{Row=6;  A=18284; B="Srta. Lara Correia";          C="Atendimento ao Cliente"; D="Problemas pessoais"; E=4; F=45078; G=10266.63}
    @{Row=7;  A=6022;  B="Dr. João Vitor Melo";         C="TI";                     D="Doença";              E=3; F=45078; G=6854.98}
    @{Row=8;  A=51531; B="Fernanda Pinto";              C="Marketing";              D="Problemas pessoais"; E=3; F=45087; G=6097.1}
    @{Row=9;  A=53029; B="Bruna Moreira";               C="Engenharia";             D="Viagem de negócios"; E=7; F=45085; G=8618.19}
    @{Row=10; A=27365; B="Ana Sophia Lima";             C="Marketing";              D="Problemas pessoais"; E=3; F=45086; G=8093.97}
    @{Row=11; A=94099; B="Emanuel Pires";               C="Engenharia";             D="Problemas pessoais"; E=4; F=45095; G=12267.98}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
